# Update market-price-derived columns (H-N) on several Leve sheets.
# Values sourced from refreshed market data snapshot (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 330.85715
$ws.Cells.Item(4, 9).Value = 378.8
$ws.Cells.Item(4, 11).Value = 378.8
$ws.Cells.Item(4, 13).Value = -264.8

$ws.Cells.Item(12, 8).Value = 200
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 14).ClearContents() | Out-Null

$ws.Cells.Item(18, 8).Value = 1007.1429
$ws.Cells.Item(18, 9).Value = 883.3333
$ws.Cells.Item(18, 10).Value = 1750
$ws.Cells.Item(18, 11).Value = 883.3333
$ws.Cells.Item(18, 12).Value = 1750
$ws.Cells.Item(18, 13).Value = -599.3333
$ws.Cells.Item(18, 14).Value = -2318

$ws.Cells.Item(33, 8).Value = 69.947365
$ws.Cells.Item(33, 9).Value = 69.947365
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 69.947365
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 159.052635
$ws.Cells.Item(33, 14).ClearContents() | Out-Null

$ws.Cells.Item(98, 8).Value = 34500
$ws.Cells.Item(98, 9).Value = 34000
$ws.Cells.Item(98, 10).Value = 35000
$ws.Cells.Item(98, 11).Value = 34000
$ws.Cells.Item(98, 12).Value = 35000
$ws.Cells.Item(98, 13).Value = -32502
$ws.Cells.Item(98, 14).Value = -37996

$ws.Cells.Item(122, 8).Value = 34500
$ws.Cells.Item(122, 9).Value = 34000
$ws.Cells.Item(122, 10).Value = 35000
$ws.Cells.Item(122, 11).Value = 102000
$ws.Cells.Item(122, 12).Value = 105000
$ws.Cells.Item(122, 13).Value = -99550
$ws.Cells.Item(122, 14).Value = -109900

$ws.Cells.Item(125, 8).Value = 500
$ws.Cells.Item(125, 9).Value = 500
$ws.Cells.Item(125, 11).Value = 4500
$ws.Cells.Item(125, 13).Value = -2040

$ws.Cells.Item(131, 8).Value = 8969
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 13).ClearContents() | Out-Null

$ws.Cells.Item(138, 8).Value = 2364.7144
$ws.Cells.Item(138, 9).Value = 754
$ws.Cells.Item(138, 11).Value = 2262
$ws.Cells.Item(138, 13).Value = 2878

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 633
$ws.Cells.Item(6, 9).Value = 599.5
$ws.Cells.Item(6, 11).Value = 599.5
$ws.Cells.Item(6, 13).Value = -426.5

$ws.Cells.Item(52, 8).Value = 19983
$ws.Cells.Item(52, 9).Value = 19975
$ws.Cells.Item(52, 10).Value = 19999
$ws.Cells.Item(52, 11).Value = 19975
$ws.Cells.Item(52, 12).Value = 19999
$ws.Cells.Item(52, 13).Value = -19657
$ws.Cells.Item(52, 14).Value = -20635

$ws.Cells.Item(88, 8).Value = 3482.5715
$ws.Cells.Item(88, 9).Value = 700
$ws.Cells.Item(88, 10).Value = 3946.3333
$ws.Cells.Item(88, 11).Value = 700
$ws.Cells.Item(88, 12).Value = 3946.3333
$ws.Cells.Item(88, 13).Value = -294
$ws.Cells.Item(88, 14).Value = -4758.3333

$ws.Cells.Item(91, 8).Value = 3482.5715
$ws.Cells.Item(91, 9).Value = 700
$ws.Cells.Item(91, 10).Value = 3946.3333
$ws.Cells.Item(91, 11).Value = 700
$ws.Cells.Item(91, 12).Value = 3946.3333
$ws.Cells.Item(91, 13).Value = 704
$ws.Cells.Item(91, 14).Value = -6754.3333

$ws.Cells.Item(97, 8).Value = 1829.5454
$ws.Cells.Item(97, 9).Value = 1315.625
$ws.Cells.Item(97, 10).Value = 3200
$ws.Cells.Item(97, 11).Value = 1315.625
$ws.Cells.Item(97, 12).Value = 3200
$ws.Cells.Item(97, 13).Value = -819.625
$ws.Cells.Item(97, 14).Value = -4192

$ws.Cells.Item(102, 8).Value = 1310.4445
$ws.Cells.Item(102, 9).Value = 1411.8572
$ws.Cells.Item(102, 10).Value = 955.5
$ws.Cells.Item(102, 11).Value = 1411.8572
$ws.Cells.Item(102, 12).Value = 955.5
$ws.Cells.Item(102, 13).Value = 210.1428000000001
$ws.Cells.Item(102, 14).Value = -4199.5

$ws.Cells.Item(132, 8).Value = 2000
$ws.Cells.Item(132, 9).Value = 2000
$ws.Cells.Item(132, 11).Value = 6000
$ws.Cells.Item(132, 13).Value = -3470

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1575.8846
$ws.Cells.Item(107, 9).Value = 1521.5238
$ws.Cells.Item(107, 10).Value = 1804.2
$ws.Cells.Item(107, 11).Value = 1521.5238
$ws.Cells.Item(107, 12).Value = 1804.2
$ws.Cells.Item(107, 13).Value = 398.4762000000001
$ws.Cells.Item(107, 14).Value = -5644.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1267
$ws.Cells.Item(16, 9).Value = 900.5
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 900.5
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = -613.5
$ws.Cells.Item(16, 14).Value = -2574

$ws.Cells.Item(33, 8).Value = 12717.8
$ws.Cells.Item(33, 9).Value = 897.25
$ws.Cells.Item(33, 11).Value = 897.25
$ws.Cells.Item(33, 13).Value = -518.25

$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).ClearContents() | Out-Null

$ws.Cells.Item(62, 8).Value = 3233.3333
$ws.Cells.Item(62, 9).Value = 3233.3333
$ws.Cells.Item(62, 11).Value = 3233.3333
$ws.Cells.Item(62, 13).Value = -2609.3333

$ws.Cells.Item(65, 8).Value = 3233.3333
$ws.Cells.Item(65, 9).Value = 3233.3333
$ws.Cells.Item(65, 11).Value = 16166.6665
$ws.Cells.Item(65, 13).Value = -13046.6665

$ws.Cells.Item(103, 8).Value = 49000
$ws.Cells.Item(103, 9).Value = 49000
$ws.Cells.Item(103, 11).Value = 49000
$ws.Cells.Item(103, 13).Value = -47828

$ws.Cells.Item(113, 8).Value = 1267
$ws.Cells.Item(113, 9).Value = 900.5
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 900.5
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = 1269.5
$ws.Cells.Item(113, 14).Value = -6340

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).ClearContents() | Out-Null
$ws.Cells.Item(3, 14).ClearContents() | Out-Null

$ws.Cells.Item(134, 8).Value = 2000
$ws.Cells.Item(134, 9).Value = 2000
$ws.Cells.Item(134, 11).Value = 6000
$ws.Cells.Item(134, 13).Value = -930

$ws.Cells.Item(140, 8).Value = 225
$ws.Cells.Item(140, 9).Value = 225
$ws.Cells.Item(140, 11).Value = 675
$ws.Cells.Item(140, 13).Value = 4505

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 10000
$ws.Cells.Item(19, 9).Value = 10000
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 13).Value = -9830

$ws.Cells.Item(22, 8).Value = 3443.1538
$ws.Cells.Item(22, 9).Value = 2599.6667
$ws.Cells.Item(22, 10).Value = 3696.2
$ws.Cells.Item(22, 11).Value = 2599.6667
$ws.Cells.Item(22, 12).Value = 3696.2
$ws.Cells.Item(22, 13).Value = -2304.6667
$ws.Cells.Item(22, 14).Value = -4286.2

$ws.Cells.Item(25, 8).Value = 10000
$ws.Cells.Item(25, 9).Value = 10000
$ws.Cells.Item(25, 11).Value = 10000
$ws.Cells.Item(25, 13).Value = -9770

$ws.Cells.Item(27, 8).Value = 3443.1538
$ws.Cells.Item(27, 9).Value = 2599.6667
$ws.Cells.Item(27, 10).Value = 3696.2
$ws.Cells.Item(27, 11).Value = 2599.6667
$ws.Cells.Item(27, 12).Value = 3696.2
$ws.Cells.Item(27, 13).Value = -2492.6667
$ws.Cells.Item(27, 14).Value = -3910.2

$ws.Cells.Item(46, 8).Value = 227109.67
$ws.Cells.Item(46, 9).Value = 1001498.5
$ws.Cells.Item(46, 11).Value = 1001498.5
$ws.Cells.Item(46, 13).Value = -1001310.5

$ws.Cells.Item(55, 8).Value = 2002
$ws.Cells.Item(55, 10).Value = 2002
$ws.Cells.Item(55, 12).Value = 2002
$ws.Cells.Item(55, 14).Value = -2348

$ws.Cells.Item(61, 8).Value = 2382.5715
$ws.Cells.Item(61, 9).Value = 2382.5715
$ws.Cells.Item(61, 11).Value = 2382.5715
$ws.Cells.Item(61, 13).Value = -2180.5715

$ws.Cells.Item(93, 8).Value = 1777.5714
$ws.Cells.Item(93, 9).Value = 2587
$ws.Cells.Item(93, 10).Value = 698.3333
$ws.Cells.Item(93, 11).Value = 2587
$ws.Cells.Item(93, 12).Value = 698.3333
$ws.Cells.Item(93, 13).Value = -1339
$ws.Cells.Item(93, 14).Value = -3194.3333

$ws.Cells.Item(113, 8).Value = 2382.5715
$ws.Cells.Item(113, 9).Value = 2382.5715
$ws.Cells.Item(113, 11).Value = 2382.5715
$ws.Cells.Item(113, 13).Value = -212.5715

$ws.Cells.Item(132, 8).Value = 2856.3333
$ws.Cells.Item(132, 9).Value = 2856.3333
$ws.Cells.Item(132, 11).Value = 8568.999899999999
$ws.Cells.Item(132, 13).Value = -6038.999899999999

$ws.Cells.Item(136, 8).Value = 4737.7144
$ws.Cells.Item(136, 9).Value = 4759.8
$ws.Cells.Item(136, 11).Value = 14279.4
$ws.Cells.Item(136, 13).Value = -11729.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 9806.166999999999
$ws.Cells.Item(113, 9).Value = 614
$ws.Cells.Item(113, 11).Value = 1842
$ws.Cells.Item(113, 13).Value = 328

Write-Output "Sheets updated via scheduled runner."
